# Applies the "quiz marksheet" update:
#  - Grading header (rows 10-12) gets real Right/Wrong/NotAttempt/Max numbers,
#    marking-scheme values, and a computed Total; row labels (No./Marking/Total)
#    pick up the title style used elsewhere in the sheet.
#  - The "-1" negative-marking cell was stored as text; store it as a real
#    number so downstream math doesn't break on it (float input handling).
#  - The third Student/Correct-Ans block (columns G:H) is removed - only two
#    blocks of answers remain.
#  - Student answers are filled in for columns A (and D for the first few
#    rows); cells get "correctStyle" (green) when the student answer matches
#    the correct answer, "incorrectStyle" (red) when it doesn't, and stay
#    "normalStyle" when left blank (not attempted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the third Student Ans / Correct Ans block (columns G:H) ---
$ws.Range("G:H").Delete()

# --- Row 10: "No." / Right / Wrong / Not Attempt / Max ---
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 28

# --- Row 11: "Marking" scheme (marks per right / penalty per wrong) ---
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

# --- Row 12: "Total" (marks earned, penalty, score/max) ---
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "29/112"

# --- First block (columns A/B): fill in the student's answers ---
function Set-Answer($ref, $value, $correct) {
    $ws.Range($ref).Value = $value
    if ($correct) {
        $ws.Range($ref).Style = "correctStyle"
    } else {
        $ws.Range($ref).Style = "incorrectStyle"
    }
}

Set-Answer "A16" "Option D" $false
Set-Answer "A18" "Option B" $true
Set-Answer "A19" "Option C" $true
Set-Answer "A21" "Option C" $true
Set-Answer "A22" "Option D" $true
Set-Answer "A25" "Option A" $true
Set-Answer "A33" "Option D" $true
Set-Answer "A39" "Option C" $false

# --- Second block (columns D/E): only rows 16-18 still carry data ---
Set-Answer "D16" "Option A" $true
Set-Answer "D17" "Option A" $false
Set-Answer "D18" "Option D" $true

# The rest of the second block (rows 19-40) is no longer used - remove it.
$ws.Range("D19:E40").Clear()
